# Thermal.xlsx update
# - Sheet1 ("Top" plane): Width/Length go from 4 cm -> 5 cm (C2, C3)
# - Sheet2 (Shield/absorption helper sheet): thickness 0.46 -> 0.5 mm (B4),
#   ambient-delta rule-of-thumb 40 -> 90 K (B12)
# - Selections updated to match where the author last clicked before saving
#   (Sheet1!C2 and Sheet2!E5), with Sheet1 remaining the active sheet/tab.
# All other cell values in the diff are formula results that recompute
# automatically from these four input changes.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: update the copper-plane Width and Length inputs ---
$ws1.Range("C2").Value = 5
$ws1.Range("C3").Value = 5

# --- Sheet2: update plating thickness and delta-T rule of thumb ---
$ws2.Range("B4").Value = 0.5
$ws2.Range("B12").Value = 90

# --- Restore/update the saved selections on each sheet ---
[void]$ws2.Activate()
[void]$ws2.Range("E5").Select()

[void]$ws1.Activate()
[void]$ws1.Range("C2").Select()
